# Weekly refresh of the Ciboulette (Hortaliza) price series:
# a new week's worth of data (rows for dates 45148) is inserted at the
# top of the existing series, pushing the older rows down by two and
# appending two "new" rows (15-16) at the bottom that re-use the data
# that used to live in rows 7-8 for the oldest still-tracked dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated values for existing rows 7-14 (shifted down by the two
#     freshly inserted observations) ---------------------------------
$rows = @{
    7  = @{ D = 45148; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 }
    8  = @{ D = 45148; I = "Segunda"; J = 60;  K = 2000; L = 2000; M = 2000; P = 667 }
    9  = @{ D = 45135; I = "Primera"; J = 70;  K = 2500; L = 2500; M = 2500; P = 833 }
    10 = @{ D = 44832; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 }
    11 = @{ D = 44832; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; P = 333 }
    12 = @{ D = 45146; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; P = 833 }
    13 = @{ D = 45146; I = "Segunda"; J = 80;  K = 2000; L = 2000; M = 2000; P = 667 }
    14 = @{ D = 44846; I = "Primera"; J = 200; K = 1200; L = 1300; M = 1250; P = 417 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}

# --- Brand-new rows 15 and 16, carrying the values that used to be in
#     (old) rows 7 and 8 before the shift -----------------------------
$newRows = @{
    15 = @{ A = 7; B = "Terminal Hortofrutícola Agro Chillán"; C = "Ñuble"; D = 44846; E = 16; F = 100112039; G = "Ciboulette"; H = "Sin especificar"; I = "Segunda"; J = 150; K = 1000; L = 1000; M = 1000; N = "`$/docena de atados"; O = "Región Metropolitana"; P = 333; Q = 3; R = "Hortaliza" }
    16 = @{ A = 7; B = "Terminal Hortofrutícola Agro Chillán"; C = "Ñuble"; D = 45133; E = 16; F = 100112039; G = "Ciboulette"; H = "Sin especificar"; I = "Primera"; J = 80;  K = 2500; L = 2500; M = 2500; N = "`$/docena de atados"; O = "Región Metropolitana"; P = 833; Q = 3; R = "Hortaliza" }
}

$dateFormat = $ws.Range("D14").NumberFormat()

foreach ($r in 15, 16) {
    $vals = $newRows[$r]
    $ws.Range("A$r").Value = $vals.A
    $ws.Range("B$r").Value = $vals.B
    $ws.Range("C$r").Value = $vals.C
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("D$r").NumberFormat = $dateFormat
    $ws.Range("E$r").Value = $vals.E
    $ws.Range("F$r").Value = $vals.F
    $ws.Range("G$r").Value = $vals.G
    $ws.Range("H$r").Value = $vals.H
    $ws.Range("I$r").Value = $vals.I
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
}

Write-Output "Updated rows 7-14 and appended rows 15-16"
